# Fix the KDA columns (A = KILLS, E = DEATHS, F = ASSISTS) for rows 2-41:
# columns A and E were stored as text ("0", "1", ...) and are converted to
# real numbers here; a couple of rows (36 and 41) additionally had the
# literal placeholder text "erro" (and row 24 had "a" in DEATHS) which are
# replaced with the correct numeric KDA values. Column F (ASSISTS) keeps
# its text storage, but the two "erro" entries become the correct numeric
# text values "10" / "11".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# KILLS (column A), numeric, rows 2..41
$killsValues = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,3,3,3,3,3)

# DEATHS (column E), numeric, rows 2..41
$deathsValues = @(0,0,0,0,0,0,0,0,0,1,1,2,2,2,2,2,2,2,2,3,3,3,4,4,4,4,5,5,6,6,6,6,6,6,6,7,7,7,7,7)

$startRow = 2
for ($i = 0; $i -lt $killsValues.Length; $i++) {
    $row = $startRow + $i

    # Column A (KILLS) -> numeric
    $ws.Cells.Item($row, 1).Value = $killsValues[$i]

    # Column E (DEATHS) -> numeric
    $ws.Cells.Item($row, 5).Value = $deathsValues[$i]
}

# Column F (ASSISTS) stays text; only the two "erro" cells change value,
# to the correct (still textual) KDA numbers.
$ws.Cells.Item(36, 6).NumberFormat = "@"
$ws.Cells.Item(36, 6).Value = "10"

$ws.Cells.Item(41, 6).NumberFormat = "@"
$ws.Cells.Item(41, 6).Value = "11"

Write-Output "KDA columns fixed"
